$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7").Value = "fgi-h"
$ws.Range("C7").Value = "nary1408"
$ws.Range("A8").Value = "fgi-p"
$ws.Range("B7").Value = "murilo121212"
$ws.Range("B8").Value = "murilo121212"
$ws.Range("C8").Value = "nary1408"
$ws.Columns("A").ColumnWidth = 9.17
$ws.Columns("B").ColumnWidth = 12
$ws.Columns("C").ColumnWidth = 8
$ws.Range("B8").Select()
